$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "fantasy points" column (E) to column G,
# and insert the new "height" and "weight" columns at E and F.

# Header row: give the two new header cells the same formatting as the
# existing header cells (bold / centered / bordered style) before touching
# their contents.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G1").Value = $ws.Range("E1").Value2
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Data rows 2-15: copy fantasy-points values from E to G, then set new values
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 255
}
